# Update the "cryptos" price/volume snapshot (GitHub Actions scheduled refresh).
# Price (D) and Volume/1h (E) columns are stored as TEXT in the workbook (e.g.
# "69.296.39", "  +2.32%  "), so numeric-looking prices are forced to text via
# NumberFormat "@" before assignment, then ClearFormats() strips the now-unneeded
# explicit format so the cell's style index matches the original (no "s" attr).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.296.39"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").Value = "3.402.99"
$ws.Range("E3").Value = "  +2.43%  "
$ws.Range("E4").Value = "  +0.09%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "587.84"
$c.ClearFormats()
$ws.Range("E5").Value = "  +1.23%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "181.34"
$c.ClearFormats()
$ws.Range("E6").Value = "  +3.92%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.599"
$c.ClearFormats()
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("E8").Value = "  +0.06%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.198"
$c.ClearFormats()
$ws.Range("E9").Value = "  +9.74%  "
$ws.Range("E10").Value = "  +2.89%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "48.59"
$c.ClearFormats()
$ws.Range("E11").Value = "  +3.92%  "
$ws.Range("E12").Value = "  +4.82%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "682.68"
$c.ClearFormats()
$ws.Range("E13").Value = "  -2.66%  "
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("D15").Value = "3.957.55"
$ws.Range("E15").Value = "  +2.54%  "
$ws.Range("D16").Value = "69.424.47"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.121"
$c.ClearFormats()
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.400.94"
$ws.Range("E18").Value = "  +2.11%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "17.74"
$c.ClearFormats()
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("E20").Value = "  +2.14%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.911"
$c.ClearFormats()
$ws.Range("E21").Value = "  +2.63%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "17.31"
$c.ClearFormats()
$ws.Range("E22").Value = "  +3.05%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.38"
$c.ClearFormats()
$ws.Range("E23").Value = "  +0.10%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "103.47"
$c.ClearFormats()
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("E26").Value = "  +2.20%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.71"
$c.ClearFormats()
$ws.Range("E27").Value = "  +3.87%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "34.03"
$c.ClearFormats()
$ws.Range("E28").Value = "  +3.25%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.80"
$c.ClearFormats()
$ws.Range("E29").Value = "  +3.59%  "
$ws.Range("E30").Value = "  +0.15%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "561.91"
$c.ClearFormats()
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("E32").Value = "  +1.81%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.63"
$c.ClearFormats()
$ws.Range("E33").Value = "  +11.90%  "
$ws.Range("E34").Value = "  +1.75%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "58.73"
$c.ClearFormats()
$ws.Range("E35").Value = "  +4.54%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").Value = "3.664.01"
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("E38").Value = "  +7.19%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "36.06"
$c.ClearFormats()
$ws.Range("E39").Value = "  +3.18%  "
$ws.Range("D40").Value = "0.0₃0727"
$ws.Range("E40").Value = "  +9.08%  "
$ws.Range("E41").Value = "  +4.57%  "
$ws.Range("E42").Value = "  +3.28%  "
$ws.Range("E43").Value = "  +2.11%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.33"
$c.ClearFormats()
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("E47").Value = "  +1.47%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.39"
$c.ClearFormats()
$ws.Range("E48").Value = "  +5.39%  "
$ws.Range("E49").Value = "  +0.06%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "132.73"
$c.ClearFormats()
$ws.Range("E50").Value = "  +1.19%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.71"
$c.ClearFormats()
$ws.Range("E51").Value = "  +3.95%  "
